$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New template-code rows appended to the master-template_type sheet
# (columns: code, descr, lang_code, is_active, cr_by, cr_dtimes)
$rows = @(
    @("RPR_DUP_UIN_EMAIL_SUB",  "Registration Failed because you have already Registered"),
    @("RPR_TEC_ISSUE_EMAIL_SUB","Re-Register because there was a Technical Issue"),
    @("RPR_UIN_REAC_EMAIL_SUB", "Uin is activated successfully"),
    @("RPR_UIN_DEAC_EMAIL_SUB", "Uin is deactivated"),
    @("RPR_UIN_GEN_EMAIL_SUB",  "UIN Generated"),
    @("RPR_UIN_UPD_EMAIL_SUB",  "UIN Details Updated")
)

$startRow = 141

# Column A (code) filled top-to-bottom first ...
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
}

# ... then column B (descr) filled bottom-to-top, matching the original
# authoring order (this also determines shared-string allocation order).
for ($i = $rows.Count - 1; $i -ge 0; $i--) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Update the window view to match the scrolled/selected state from the edit session
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F146").Select()

# Page setup changes captured with the edit
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
